# Applies the edits described by the commit "Add files via upload":
#  - Fix the TOPSIS ideal-best / ideal-worst rows: columns I, J, K (the
#    "non beneficial" criteria) had MAX/MIN swapped between row 101 (max)
#    and row 102 (min); swap the formulas back the right way round.
#  - Un-share K40 into an explicit "=K13*K13" formula.
#  - Fix the closeness-coefficient formula in column L (rows 107:124):
#    was J/H+J, should be J/(H+J).
#  - Fix the RANK.EQ formula in column M (rows 107:124) to pass an
#    explicit ascending sort order (3rd argument = 1).
#  - Resize column G and update the active selection / scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ideal best (row 101) / ideal worst (row 102): swap I, J, K ---
$ws.Range("I101").Formula = "=MIN(I83:I100)"
$ws.Range("J101").Formula = "=MIN(J83:J100)"
$ws.Range("K101").Formula = "=MIN(K83:K100)"

$ws.Range("I102").Formula = "=MAX(I83:I100)"
$ws.Range("J102").Formula = "=MAX(J83:J100)"
$ws.Range("K102").Formula = "=MAX(K83:K100)"

# --- K40: make it an explicit (non-shared) formula ---
$ws.Range("K40").Formula = "=K13*K13"

# --- L107:L124 and M107:M124 ---
For ($r = 107; $r -le 124; $r++) {
    $ws.Range("L$r").Formula = "=J$r/(H$r+J$r)"
    $ws.Range("M$r").Formula = "=RANK.EQ(L$r,L`$107:L`$124,1)+COUNTIF(L`$107:L`$124,L$r)-1"
}

# --- Column widths: split the G:H group so G gets its own width ---
$ws.Columns("G").ColumnWidth = 11.666666666666666

# --- Selection / scroll position ---
$ws.Range("O118").Select()
